$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 11.45740717551576
# Row 3
$ws.Range("H3").Value = 10.08421220545006
# Row 4
$ws.Range("H4").Value = 8.474440037169567
# Row 5
$ws.Range("H5").Value = 8.469079889133782
# Row 6
$ws.Range("H6").Value = 7.167398355129854
# Row 7
$ws.Range("D7").Value = 16
$ws.Range("E7").Value = '60863a15760523386e761cfb'
$ws.Range("F7").Value = 'Roshni'
$ws.Range("G7").Value = 'female'
$ws.Range("H7").Value = 6.194924391488313
# Row 8
$ws.Range("D8").Value = 13
$ws.Range("E8").Value = '5697d4ae7183b8000d0fc201'
$ws.Range("F8").Value = 'Tu'
$ws.Range("G8").Value = 'male'
$ws.Range("H8").Value = 6.053014121671316
# Row 9
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = '5c27de12a2b00a00018b2c16'
$ws.Range("F9").Value = 'Ankai'
$ws.Range("G9").Value = 'male'
$ws.Range("H9").Value = 5.405704693608066
# Row 10
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = '60bd88b8fc436774352f53b9'
$ws.Range("F10").Value = 'Annes'
$ws.Range("G10").Value = 'female'
$ws.Range("H10").Value = 5.190617474046819
# Row 11
$ws.Range("H11").Value = 4.0828279198773
# Row 12
$ws.Range("H12").Value = 2.154970194355574
# Row 13
$ws.Range("H13").Value = 0.2602112850569155
# Row 14
$ws.Range("D14").Value = 7
$ws.Range("E14").Value = '6024c18b094ac71dd93f4f5a'
$ws.Range("F14").Value = 'Katherine'
$ws.Range("H14").Value = 8.316648944792245
# Row 15
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = '60778ed0fde3e9c3a96f1d11'
$ws.Range("F15").Value = 'Melissa'
$ws.Range("H15").Value = 8.023344841524992
# Row 16
$ws.Range("H16").Value = 7.429443214079729
# Row 17
$ws.Range("H17").Value = 7.094459853851288
# Row 18
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = '5eeaa065c7acf61c4322f6d9'
$ws.Range("F18").Value = 'Yonifredy'
$ws.Range("H18").Value = 6.304187637973969
# Row 19
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = '5f5ea8227fa75676f56f9276'
$ws.Range("F19").Value = 'Carlos'
$ws.Range("H19").Value = 6.19077550199683
# Row 20
$ws.Range("H20").Value = 5.065205973220809
# Row 21
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = '5e0adc8f4cac6834756db412'
$ws.Range("F21").Value = 'Mary'
$ws.Range("H21").Value = 3.358111939047832
# Row 22
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = '5e706891c396cc64388ef760'
$ws.Range("F22").Value = 'Maria'
$ws.Range("H22").Value = 3.344821734808749
# Row 23
$ws.Range("H23").Value = 2.475630392065158
# Row 24
$ws.Range("H24").Value = 1.281029176420817
# Row 25
$ws.Range("H25").Value = 0.379747773547242
